# Update the "K" (strikeouts) column (G) on the active sheet with newly
# regenerated values (replacing the previous "Strike#" based figures).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    11 = 0
    13 = 2
    14 = 1
    15 = 1
    16 = 2
    17 = 0
    19 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
